$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 333.35
$ws.Range("F2").Value = 1.1
$ws.Range("N2").Value = 53.71147335634279

# Row 3
$ws.Range("D3").Value = 78.17
$ws.Range("E3").Value = 43.7
$ws.Range("F3").Value = 2.4
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 36
$ws.Range("K3").Value = 58.1
$ws.Range("N3").Value = 53.71147335634279

# Row 4
$ws.Range("D4").Value = 77.75
$ws.Range("E4").Value = 46.2
$ws.Range("F4").Value = 2.25
$ws.Range("N4").Value = 53.71147335634279

# Row 5
$ws.Range("D5").Value = 110.95
$ws.Range("E5").Value = 65.40000000000001
$ws.Range("F5").Value = 2.76
$ws.Range("G5").Value = 50
$ws.Range("K5").Value = 47.1
$ws.Range("N5").Value = 53.71147335634279
